$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 275 (shifts the existing data, incl. row 333,
# down by one row) and populate it with the new Coimbatore (CJB) colo entry.
$ws.Rows.Item(275).Insert()

$ws.Cells.Item(275, 1).Value = "CJB"
$ws.Cells.Item(275, 2).Value = "Coimbatore, India"
$ws.Cells.Item(275, 3).Value = "Asia"
$ws.Cells.Item(275, 4).Value = "Coimbatore"
$ws.Cells.Item(275, 5).Value = "India"
$ws.Cells.Item(275, 6).Value = "IN"

# Match the formatting used by column A for the other data rows (bordered /
# bold / centered "colo" column style).
$colo = $ws.Cells.Item(275, 1)
$colo.Font.Bold = $true
$colo.HorizontalAlignment = -4108
$colo.VerticalAlignment = -4160
$colo.Borders.LineStyle = 1
